$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (final save differential) values per repulled data
$ws.Range("F6").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("F11").Value = 2
$ws.Range("F15").Value = -1
$ws.Range("F18").Value = 3
$ws.Range("F20").Value = 3
$ws.Range("F31").Value = 1
$ws.Range("F32").Value = -1
$ws.Range("F34").Value = 2
$ws.Range("F40").Value = 0
$ws.Range("F42").Value = 3
$ws.Range("F47").Value = 0
$ws.Range("F48").Value = 4
